# Auto-generated edit script for cryptos.xlsx update
# Applies the cell-level text changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character (U+2083) used in the PEPE price cell (D28).
$sub3 = [string][char]0x2083

$ws.Range("D2").Value = '59.745.14'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '2.403.24'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.37'
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.89'
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +4.04%  '
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.71'
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("E13").Value = '  +2.12%  '
$ws.Range("D14").Value = '2.829.78'
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '59.674.97'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").Value = '2.405.98'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.31'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '328.30'
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("E21").Value = '  -4.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.51'
$ws.Range("E23").Value = '  +2.99%  '
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.65'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").Value = '0.0' + $sub3 + '0770'
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.29'
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.07'
$ws.Range("E31").Value = '  -3.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.59'
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("E33").Value = '  -1.53%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -1.48%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.60'
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '313.81'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("E40").Value = '  -3.34%  '
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.97'
$ws.Range("E42").Value = '  -2.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0970'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.53'
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0514'
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.576'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0222'
$ws.Range("E47").Value = '  -1.82%  '
$ws.Range("E48").Value = '  -4.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.52'
$ws.Range("E49").Value = '  -2.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.04'
$ws.Range("E50").Value = '  -0.13%  '
$ws.Range("E51").Value = '  -3.35%  '
